{"js": "// Replace the outdated notary instructions sentence with the updated\n// guidance, and drop the stray \"_GoBack\" bookmark left over from the\n// previous edit session (QA testing feedback).\n\nconst body = context.document.body;\n\n// The sentence previously read:\n//   \"Your notary has to see you sign it, and they must sign it right\n//    after you do.\"\n// and must become:\n//   \"Your notary must see you sign the form, then sign it after you.\n//    Do not sign your forms before they are notarized.\"\nconst oldText =\n  \"has to see you sign it, and they must sign it right after you do.\";\nconst newText =\n  \"must see you sign the form, then sign it after you. \" +\n  \"Do not sign your forms before they are notarized.\";\n\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the notary instructions sentence to update.\");\n}\n\nresults.items[0].insertText(newText, Word.InsertLocation.replace);\n\n// Remove the leftover \"_GoBack\" bookmark.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Replace the outdated notary instructions sentence with the updated\n# guidance, and drop the stray \"_GoBack\" bookmark left over from the\n# previous edit session (QA testing feedback).\n\n$d = $word.ActiveDocument\n\n$oldText = \"has to see you sign it, and they must sign it right after you do.\"\n$newText = \"must see you sign the form, then sign it after you. Do not sign your forms before they are notarized.\"\n\n$find = $d.Content.Find\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
